$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $r = $ws.Range($range)
    $r.NumberFormat = "@"
    $r.Value = $value
}

# Row 2 - BNB
Set-TextValue "D2" "246.45"

# Row 3 - OKB
Set-TextValue "D3" "21.90"

# Row 4 - HuobiToken
Set-TextValue "D4" "5.404"

# Row 5 - Cronos
Set-TextValue "D5" "0.05783"

# Row 6
Set-TextValue "D6" "3.386"

# Row 8
Set-TextValue "D8" "0.8087"

# Row 9 - FTXToken
Set-TextValue "D9" "0.9539"
$ws.Range("E9").Value = "8FTXTokenFTTBestin24h"

# Row 11
Set-TextValue "D11" "0.07486"

# Row 12
Set-TextValue "D12" "0.03188"

# Row 13
Set-TextValue "D13" "0.03013"

# Row 14
Set-TextValue "D14" "4.164"

# Row 15
Set-TextValue "D15" "0.09410"

# Row 16
Set-TextValue "D16" "0.001586"

# Row 17
Set-TextValue "D17" "0.04813"

# Row 18
Set-TextValue "D18" "0.0005898"

# Row 19
Set-TextValue "D19" "0.006182"

# Row 20
Set-TextValue "D20" "0.004112"

# Row 21
Set-TextValue "D21" "0.0009990"

# Row 23
Set-TextValue "D23" "3.770"

# Row 24
Set-TextValue "D24" "2.231"

# Row 26
Set-TextValue "D26" "0.1259"

# Row 27 - UpBots
Set-TextValue "D27" "0.0003032"
$ws.Range("E27").Value = "26UpBotsUBXT"

# Row 40
Set-TextValue "D40" "0.03882"

# Row 41 - becomes BKEXToken
$ws.Range("B41").Value = "BKEXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue "D41" "0.1075"
$ws.Range("E41").Value = "40BKEXTokenBKK"

# Row 42 - becomes CEJI
$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue "D42" "0.002620"
$ws.Range("E42").Value = "41CEJICEJI"

# Row 43 - becomes KickToken
$ws.Range("B43").Value = "KickToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
Set-TextValue "D43" "0.003031"
$ws.Range("E43").Value = "42KickTokenKICKWorstin24h"

# Row 44
Set-TextValue "D44" "0.006195"

# Row 45
Set-TextValue "D45" "0.00005586"

# Row 49
Set-TextValue "D49" "0.00002099"
